$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Fill in the Actual output and Result for row 11 (SL. No 10)
$ws.Range("F11").Value = "It gets displayed as expected"
$ws.Range("G11").Value = "Pass"

# Update the saved view state (top-left cell and selection)
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F12").Select()
